# Selenium_AutomaçãoDeTestes.docx — apply commit "Testes Selenium IF ELSE"
$d = $word.ActiveDocument

# --- wdReplaceOne constant used by Find.Execute(Replace:=) ---
$wdReplaceOne = 1

# 1) "Pirâmide de Escopo" -> "Pirâmide Escopo"
$d.Content.Find.Execute("Pirâmide de Escopo", $false, $false, $false, $false, $false, $true, 1, $false, "Pirâmide Escopo", $wdReplaceOne) | Out-Null

# 2) "Validação da Tela (interface) de testes:" -> "Validação da Tela de testes:"
$d.Content.Find.Execute("Validação da Tela (interface) de testes:", $false, $false, $false, $false, $false, $true, 1, $false, "Validação da Tela de testes:", $wdReplaceOne) | Out-Null

# 3) "...campos principais da interface web..." -> "...campos principais do formulário web..."
$d.Content.Find.Execute("campos principais da interface web", $false, $false, $false, $false, $false, $true, 1, $false, "campos principais do formulário web", $wdReplaceOne) | Out-Null

# 4) "store: ARMAZENAMENTO;" -> "store: ARMAZENAMENTO: armazena o valor de um elemento em uma variável;"
$d.Content.Find.Execute(": ARMAZENAMENTO;", $false, $false, $false, $false, $false, $true, 1, $false, ": ARMAZENAMENTO: armazena o valor de um elemento em uma variável;", $wdReplaceOne) | Out-Null

# 5) "...independentemente do valor que o elemento em questão possua;" -> "...independentemente do tipo ou do valor do elemento em questão;"
$d.Content.Find.Execute("independentemente do valor que o elemento em questão possua;", $false, $false, $false, $false, $false, $true, 1, $false, "independentemente do tipo ou do valor do elemento em questão;", $wdReplaceOne) | Out-Null

# 6) "assertText: valida a presença de textos." -> "assertText: validação da presença de textos em campos."
$d.Content.Find.Execute("valida a presença de textos.", $false, $false, $false, $false, $false, $true, 1, $false, "validação da presença de textos em campos.", $wdReplaceOne) | Out-Null

# --- Append new paragraphs after the current last paragraph ---

function Add-PlainParagraph([string]$text) {
    $lastP = $d.Paragraphs.Last
    $lastP.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Last
    $newP.Range.Text = $text
    return $newP
}

function Add-EmptyParagraph() {
    $lastP = $d.Paragraphs.Last
    $lastP.Range.InsertParagraphAfter()
    return $d.Paragraphs.Last
}

function Add-BoldParagraph([string]$text) {
    $lastP = $d.Paragraphs.Last
    $lastP.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Last
    $newP.Range.Text = $text
    $newP.Range.Bold = 1
    $newP.Range.BoldBi = 1
    return $newP
}

# empty spacer paragraph
Add-EmptyParagraph | Out-Null

# "Scripts de testes são..."
Add-PlainParagraph "Scripts de testes são o passo-a-passo dos testes gravados no Selenium. Podemos gravar diversos scripts de diversos testes." | Out-Null

# "Boas práticas:"
Add-PlainParagraph "Boas práticas:" | Out-Null

# bullet list of good practices
Add-PlainParagraph "- Validação da tela de testes antes de iniciar um teste;" | Out-Null
Add-PlainParagraph "- Testar os fluxos principais básicos primeiramente e, depois, proceder aos fluxos alternativos diversos que estejam no escopo do projeto;" | Out-Null

# "- Nomear os casos de testes seguindo o padrão CT01NomeDoCasoDeTeste;" (CT01 bold+underline, NomeDoCasoDeTeste underline)
$p = Add-PlainParagraph "- Nomear os casos de testes seguindo o padrão CT01NomeDoCasoDeTeste;"
$pStart = $p.Range.Start
$prefix = "- Nomear os casos de testes seguindo o padrão "
$ct01Start = $pStart + $prefix.Length
$ct01End = $ct01Start + [string]"CT01".Length
$nomeEnd = $ct01End + [string]"NomeDoCasoDeTeste".Length
$ct01Range = $d.Range($ct01Start, $ct01End)
$ct01Range.Bold = 1
$ct01Range.BoldBi = 1
$ct01Range.Underline = 1
$nomeRange = $d.Range($ct01End, $nomeEnd)
$nomeRange.Underline = 1

Add-PlainParagraph "- Comentar ao máximo possível as ações dos scripts de testes, a fim de facilitar a utilização e a manutenção posterior dos desses." | Out-Null

Add-PlainParagraph "Para inserir comentários no script de testes no Selenium IDE 3.17.0, rodando no navegador Edge, basta inserir uma linha de comando iniciar a descrição com duas barras “//” no campo “Command”." | Out-Null

# empty spacer paragraph
Add-EmptyParagraph | Out-Null

# "Plugins" heading (bold)
Add-BoldParagraph "Plugins" | Out-Null

Add-PlainParagraph "O Selenium IDE possui uma variedade de plugins para browsers ou arquivos .js, cujo caminho é apontado no menu “Opções” da IDE. Abaixo, podemos verificar alguns desses plugins:" | Out-Null

Add-PlainParagraph "- Selenium IDE Button: Permite alternar a exibição da IDE do Selenium, seja em um pop-up ou no próprio frame do browser;" | Out-Null
Add-PlainParagraph "- Flow Control: Adiciona comandos de repetição ao script de teste;" | Out-Null
Add-PlainParagraph "- ScreenShot on Fail: Registra um print da tela quando ocorre um erro na execução do teste;" | Out-Null
Add-PlainParagraph "- Pretty Report: Exporta os resultados de testes em um relatório com um visual mais bonito e legível." | Out-Null

# "Obs.: Aula de 2021..." ("Obs.: " bold, rest normal)
$p2 = Add-PlainParagraph "Obs.: Aula de 2021, portanto, pode ser que alguns plugins não existam ou não sejam mais necessários."
$p2Start = $p2.Range.Start
$obsLen = [string]"Obs.: ".Length
$obsRange = $d.Range($p2Start, $p2Start + $obsLen)
$obsRange.Bold = 1
$obsRange.BoldBi = 1

# empty spacer paragraph
Add-EmptyParagraph | Out-Null

# "Estruturas Condicionais" heading (bold)
Add-BoldParagraph "Estruturas Condicionais" | Out-Null

# last paragraph with "if-else" in italics
$p3 = Add-PlainParagraph "O Selenium IDE 3.17.0, rodando em Edge, nos possibilita criar estruturas condicionais com os elementos da tela. Por exemplo, se quisermos verificar se um campo possui ou não informação, podemos fazer uma if-else para verificar essa condição. Diversas outras aplicações de condição podem ser efetuadas."
$p3Start = $p3.Range.Start
$prefix3 = "O Selenium IDE 3.17.0, rodando em Edge, nos possibilita criar estruturas condicionais com os elementos da tela. Por exemplo, se quisermos verificar se um campo possui ou não informação, podemos fazer uma "
$ifStart = $p3Start + $prefix3.Length
$ifEnd = $ifStart + [string]"if-else".Length
$ifRange = $d.Range($ifStart, $ifEnd)
$ifRange.Italic = 1

Write-Output "All edits applied successfully."
